$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.332.14"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.72"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.41"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("E7").Value = "  -2.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4027"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.63"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07841"
$ws.Range("E10").Value = "  -2.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9819"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.31"
$ws.Range("E12").Value = "  -1.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.09"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.832"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.004"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.20"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06548"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001018"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.17"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9984"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.326.57"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.332"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.83"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.249"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.085.90"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.34"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.28"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.061"
$ws.Range("E29").Value = "  -4.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.277"
$ws.Range("E30").Value = "  -4.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.19"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9531"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09322"
$ws.Range("E33").Value = "  -1.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.597"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.384"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.213"
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06014"
$ws.Range("E37").Value = "  -1.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02204"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.278"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5742"
$ws.Range("E42").Value = "  -3.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1805"
$ws.Range("E43").Value = "  -3.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.01"
$ws.Range("E44").Value = "  -3.79%  "
$ws.Range("E45").Value = "  -4.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.293"
$ws.Range("E46").Value = "  +13.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5414"
$ws.Range("E47").Value = "  -3.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.84"
$ws.Range("E48").Value = "  -3.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07166"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.879"
$ws.Range("E50").Value = "  -4.56%  "
$ws.Range("E51").Value = "  -0.79%  "
